$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (C and D) for EPRetPrem / IntRetPrem,
# shifting MktPremium/ProdRet/ResRet/RetWells right by two columns.
$ws.Columns("C:D").Insert()

# New header cells
$ws.Range("C1").Value = "EPRetPrem"
$ws.Range("D1").Value = "IntRetPrem"

# New data values (matlab commands for clipboard)
$ws.Range("C2").Value = -0.40672404853143224
$ws.Range("D2").Value = -1990.9054517596526
$ws.Range("C3").Value = -0.12331656348690673
$ws.Range("D3").Value = -1992.0237147063688
$ws.Range("C4").Value = -0.07233640452835807
$ws.Range("D4").Value = -1992.862124027988
$ws.Range("C5").Value = -0.28112312216259605
$ws.Range("D5").Value = -1993.9851908843464
$ws.Range("C6").Value = 0.09541476065886806
$ws.Range("D6").Value = -1994.774299719904
$ws.Range("C7").Value = 0.22092369660558633
$ws.Range("D7").Value = -1995.818303705138
$ws.Range("C8").Value = -0.14768340458183332
$ws.Range("D8").Value = -1996.8175564642638
$ws.Range("C9").Value = -0.4440882555824529
$ws.Range("D9").Value = -1997.9429986830105
$ws.Range("C10").Value = 0.11300284903673877
$ws.Range("D10").Value = -1998.8622551868068
$ws.Range("C11").Value = 0.39226784987588126
$ws.Range("D11").Value = -1999.9244003603874
$ws.Range("C12").Value = -0.28086053266448835
$ws.Range("D12").Value = -2001.0832951564967
$ws.Range("C13").Value = -0.04359004773124067
$ws.Range("D13").Value = -2002.1586287406246
$ws.Range("C14").Value = 0.18959209869171162
$ws.Range("D14").Value = -2002.7942619943894
$ws.Range("C15").Value = 0.2743596446669518
$ws.Range("D15").Value = -2003.7723046200983
$ws.Range("C16").Value = 0.47076701936886745
$ws.Range("D16").Value = -2004.859592737125
$ws.Range("C17").Value = -0.008264917599127598
$ws.Range("D17").Value = -2005.72286975985
$ws.Range("C18").Value = 0.31700140929592346
$ws.Range("D18").Value = -2006.7583219905407
$ws.Range("C19").Value = -0.44453300322935135
$ws.Range("D19").Value = -2008.2676692398031
$ws.Range("C20").Value = 0.3412242238375859
$ws.Range("D20").Value = -2009.0405573324622
$ws.Range("C21").Value = 0.08050396277358625
$ws.Range("D21").Value = -2009.8567614417334
$ws.Range("C22").Value = -0.07318155143480494
$ws.Range("D22").Value = -2010.8882282663687
$ws.Range("C23").Value = 0.019752696199697113
$ws.Range("D23").Value = -2012.005618743879
$ws.Range("C24").Value = 0.22761349649177387
$ws.Range("D24").Value = -2012.8333337165604
$ws.Range("C25").Value = -0.1272678823388888
$ws.Range("D25").Value = -2014.099589207101

# Restore/update the view selection to match the edited range
$ws.Range("A1:H25").Select()
